# The edit targets the "NearestNeighbors" worksheet (4th tab, the
# tab that was already active/selected in the workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NearestNeighbors")
$ws.Activate()

# Update the current selection from the stale "L20" to the data range
# A2:S10, anchored with A2 as the active cell - matches the diff's
# <selection activeCell="A2" sqref="A2:S10"/>.
$ws.Range("A2:S10").Select()

# The RANDBETWEEN(1,1000) formulas scattered across the table are
# volatile, so simply forcing a recalculation refreshes their cached
# <v> values (same effect as the author re-saving the workbook after
# Excel recalculated the random numbers).
$excel.CalculateFull()
